$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 9
$ws.Range("B3").Value = 1481
$ws.Range("C3").Value = 3328
$ws.Range("D3").Value = 16805

# Row 4
$ws.Range("A4").Value = 283.64999999999998
$ws.Range("B4").Value = 1858.5
$ws.Range("C4").Value = 4581.8
$ws.Range("D4").Value = 18369.099999999999

# Row 6
$ws.Range("A6").Value = 50.03
$ws.Range("B6").Value = 72.19
$ws.Range("C6").Value = 73.45
$ws.Range("D6").Value = 72.59

# Row 7
$ws.Range("A7").Value = 51.15
$ws.Range("B7").Value = 74.3
$ws.Range("C7").Value = 64.19
$ws.Range("D7").Value = 64.62

$excel.CalculateFullRebuild()
